$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 498 (pushes old rows 498:537 down to 499:538,
# dimension grows from A1:R537 to A1:R538).
$ws.Rows.Item(498).Insert()

# Populate the newly inserted row with the new weekly price observation.
$ws.Cells.Item(498, 1).Value  = 3
$ws.Cells.Item(498, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(498, 3).Value  = 'Coquimbo'
$ws.Cells.Item(498, 4).Value  = 45013
$ws.Cells.Item(498, 5).Value  = 5
$ws.Cells.Item(498, 6).Value  = 100112040
$ws.Cells.Item(498, 7).Value  = 'Cilantro'
$ws.Cells.Item(498, 8).Value  = 'Sin especificar'
$ws.Cells.Item(498, 9).Value  = 'Primera'
$ws.Cells.Item(498, 10).Value = 170
$ws.Cells.Item(498, 11).Value = 4500
$ws.Cells.Item(498, 12).Value = 5000
$ws.Cells.Item(498, 13).Value = 4824
$ws.Cells.Item(498, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(498, 15).Value = 'La Cruz'
$ws.Cells.Item(498, 16).Value = 1608
$ws.Cells.Item(498, 17).Value = 3
$ws.Cells.Item(498, 18).Value = 'Hortaliza'
